$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.689.32'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.584.86'
$ws.Range("E3").Value = '  -1.99%  '
$ws.Range("E4").Value = '  +1.36%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.11'
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("E6").Value = '  -1.96%  '
$ws.Range("E7").Value = '  +1.43%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.20'
$ws.Range("E8").Value = '  -3.71%  '
$ws.Range("E9").Value = '  -0.66%  '
$ws.Range("E10").Value = '  -2.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0866'
$ws.Range("E11").Value = '  -0.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.809.83'
$ws.Range("E12").Value = '  -1.95%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.586.77'
$ws.Range("E13").Value = '  -1.87%  '
$ws.Range("E14").Value = '  -2.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.531'
$ws.Range("E15").Value = '  -4.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.653.64'
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("E17").Value = '  -2.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '219.48'
$ws.Range("E18").Value = '  -3.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0692'
$ws.Range("E19").Value = '  -3.13%  '
$ws.Range("E20").Value = '  -4.47%  '
$ws.Range("E21").Value = '  +1.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.14'
$ws.Range("E22").Value = '  -3.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.54'
$ws.Range("E23").Value = '  -4.99%  '
$ws.Range("E24").Value = '  -2.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.90'
$ws.Range("E25").Value = '  +0.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.82'
$ws.Range("E26").Value = '  -0.86%  '
$ws.Range("E27").Value = '  +1.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.10'
$ws.Range("E28").Value = '  -2.07%  '
$ws.Range("E30").Value = '  -1.40%  '
$ws.Range("E31").Value = '  -2.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.25'
$ws.Range("E32").Value = '  -3.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.381.30'
$ws.Range("E33").Value = '  -0.57%  '
$ws.Range("E34").Value = '  -4.40%  '
$ws.Range("E35").Value = '  -3.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.979'
$ws.Range("E36").Value = '  -1.56%  '
$ws.Range("E37").Value = '  +0.19%  '
$ws.Range("E38").Value = '  -2.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.538'
$ws.Range("E39").Value = '  -2.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.822'
$ws.Range("E40").Value = '  -2.22%  '
$ws.Range("E41").Value = '  +1.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.976'
$ws.Range("E42").Value = '  -3.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.77'
$ws.Range("E43").Value = '  -2.33%  '
$ws.Range("E44").Value = '  +2.76%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.43'
$ws.Range("E45").Value = '  -3.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.22'
$ws.Range("E46").Value = '  -2.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.720.80'
$ws.Range("E47").Value = '  -1.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.39'
$ws.Range("E48").Value = '  +0.90%  '
$ws.Range("E49").Value = '  +11.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0971'
$ws.Range("E50").Value = '  -3.57%  '
$ws.Range("E51").Value = '  -0.54%  '
